# Updated the guide to datasets.xlsx
# - Replaces the old "GDP per capita" (row 16) and "Population growth" (row 17)
#   dataset entries with newer datasets (Median earnings / Household projections)
# - Appends seven new dataset rows (18-24): existing social housing supply,
#   five years of industry employment data, and an over-65s population dataset
# - All new "Access date" entries use 12/06/2024 (serial 45455)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DatasetRow($Row, $Dataset, $Usage, $AccessDateSerial, $FileName, $Hyperlink) {
    $ws.Cells.Item($Row, 1).Value = $Dataset
    $ws.Cells.Item($Row, 2).Value = $Usage
    $ws.Cells.Item($Row, 4).Value = $FileName
    $ws.Cells.Item($Row, 5).Value = $Hyperlink

    # Copy the number format from an existing formatted date cell (row 2, col C)
    # so the new cell re-uses the workbook's existing date style rather than
    # creating a brand-new custom number format.
    $ws.Cells.Item(2, 3).Copy() | Out-Null
    $dateCell = $ws.Cells.Item($Row, 3)
    $dateCell.PasteSpecial(-4122) | Out-Null
    $dateCell.Value = $AccessDateSerial
}

# Row 16: was "GDP per capita" / ONS regional GDP dataset -> now Median earnings
Set-DatasetRow 16 "ONS Median earnings" "Median earnings" 45455 "median_earnings.csv" "https://www.ons.gov.uk/peoplepopulationandcommunity/housing/datasets/ratioofhousepricetoresidencebasedearningslowerquartileandmedian/current"

# Row 17: was "Population growth" / ONS population estimates -> now Household projections
Set-DatasetRow 17 "ONS Household projections for England: 2018-based principal projection edition" "Households total / Households change" 45455 "2018basedhhpsprincipalprojection.xlsx" "https://www.ons.gov.uk/peoplepopulationandcommunity/populationandmigration/populationprojections/datasets/householdprojectionsforengland"

# Row 18: Subnational estimates of dwellings by tenure
Set-DatasetRow 18 "Subnational estimates of dwellings by tenure, England, 2012 to 2021" "Existing social housing supply" 45455 "subnationaldwellingsbytenure2021.xlsx" "https://www.ons.gov.uk/peoplepopulationandcommunity/housing/datasets/subnationaldwellingstockbytenureestimates"

# Rows 19-23: Business Register and Employment Survey, one row per year 2016-2020
Set-DatasetRow 19 "2016 Business Register and Employment Survey: Broad industrial group" "Professional and financial employment" 45455 "2016_industry_employment.csv" "Nomis data query"
Set-DatasetRow 20 "2017 Business Register and Employment Survey: Broad industrial group" "Professional and financial employment" 45455 "2017_industry_employment.csv" "Nomis data query"
Set-DatasetRow 21 "2018 Business Register and Employment Survey: Broad industrial group" "Professional and financial employment" 45455 "2018_industry_employment.csv" "Nomis data query"
Set-DatasetRow 22 "2019 Business Register and Employment Survey: Broad industrial group" "Professional and financial employment" 45455 "2019_industry_employment.csv" "Nomis data query"
Set-DatasetRow 23 "2020 Business Register and Employment Survey: Broad industrial group" "Professional and financial employment" 45455 "2020_industry_employment.csv" "Nomis data query"

# Row 24: ONS population estimates by single year of age (over 65s)
Set-DatasetRow 24 "ONS: Population estimates - local authority based by single year of age" "Over 65s percentage" 45455 "la_all_ages.csv" "Nomis data query"

# Move the active selection to the last row added, matching the saved view state
$ws.Range("A24").Select() | Out-Null
